# Fix some logics in thongKe - Fix glitch display
# Update the invoice number and date/time text on the "HoaDonMau" sheet
# to reflect the corrected values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HoaDonMau")

$ws.Range("A6").Value = "Số HD: HD121223001"
$ws.Range("A7").Value = "Ngày giờ: 23:08:23 - 12/12/2023"
